$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B / C updates (row 31 / 32 swap: Fetch.AI <-> Aptos) ---
$ws.Range("B31").Value = "Aptos"
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"

# --- Column D (Price) updates: force text format so values like "1.20" / "58.375.20" are not coerced to numbers ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.375.20"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.488.85"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "521.95"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.65"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.508.83"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.341"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.932.43"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.304.13"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.20"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.497.25"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.74"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "322.65"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.76"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.414"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.43"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0756"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.39"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.20"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.70"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.17"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.33"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.05"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.66"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.806"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.22"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "278.21"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.48"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.600"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "125.02"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.30"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.746.99"

# --- Column E (Volume 1h) updates ---
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  -1.97%  "
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("E8").Value = "  -1.17%  "
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("E10").Value = "  -1.67%  "
$ws.Range("E11").Value = "  -0.76%  "
$ws.Range("E12").Value = "  -1.44%  "
$ws.Range("E13").Value = "  -1.56%  "
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("E15").Value = "  -0.84%  "
$ws.Range("E16").Value = "  -2.19%  "
$ws.Range("E17").Value = "  -1.11%  "
$ws.Range("E18").Value = "  -0.88%  "
$ws.Range("E19").Value = "  -2.48%  "
$ws.Range("E20").Value = "  -0.74%  "
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("E23").Value = "  -2.35%  "
$ws.Range("E24").Value = "  -0.34%  "
$ws.Range("E25").Value = "  -0.94%  "
$ws.Range("E26").Value = "  -0.90%  "
$ws.Range("E27").Value = "  -0.54%  "
$ws.Range("E28").Value = "  -0.88%  "
$ws.Range("E29").Value = "  -0.44%  "
$ws.Range("E30").Value = "  -0.24%  "
$ws.Range("E31").Value = "  -0.82%  "
$ws.Range("E32").Value = "  +6.02%  "
$ws.Range("E33").Value = "  -1.79%  "
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("E36").Value = "  -0.76%  "
$ws.Range("E37").Value = "  -1.57%  "
$ws.Range("E38").Value = "  +0.43%  "
$ws.Range("E39").Value = "  -2.00%  "
$ws.Range("E40").Value = "  -0.30%  "
$ws.Range("E41").Value = "  +1.43%  "
$ws.Range("E42").Value = "  +5.16%  "
$ws.Range("E43").Value = "  -1.20%  "
$ws.Range("E44").Value = "  -2.37%  "
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("E47").Value = "  -1.19%  "
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("E49").Value = "  -1.07%  "
$ws.Range("E50").Value = "  +0.68%  "
$ws.Range("E51").Value = "  -0.02%  "
